# Rename the sheet and rewrite its data per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Shahbaz Ahmed"

$headers = @("matchNo","teamName","batterName","states","runs","balls","fours","sixes","sr","opponentTeamName","venue","date","result")

$data = @(
    @("Eliminator","Royal Challengers Bangalore","Shahbaz Ahmed","c Shivam Mavi b Ferguson","13","14","1","0","92.85","Kolkata Knight Riders","Sharjah","October 11","KKR won by 4 wickets (with 2 balls remaining)"),
    @("52nd","Royal Challengers Bangalore","Shahbaz Ahmed","c Williamson b Holder","14","9","2","0","155.55","Sunrisers Hyderabad","Abu Dhabi","October 06","Sunrisers won by 4 runs"),
    @("48th","Royal Challengers Bangalore","Shahbaz Ahmed","b Mohammed Shami","8","4","0","1","200.00","Punjab Kings","Sharjah","October 03","RCB won by 6 runs"),
    @("26th","Royal Challengers Bangalore","Shahbaz Ahmed","c Harpreet Brar b Ravi Bishnoi","8","11","1","0","72.72","Punjab Kings","Ahmedabad","April 30","Punjab Kings won by 34 runs"),
    @("6th","Royal Challengers Bangalore","Shahbaz Ahmed","c Rashid Khan b Nadeem","14","10","0","1","140.00","Sunrisers Hyderabad","Chennai","April 14","RCB won by 6 runs"),
    @("39th","Royal Challengers Bangalore","Shahbaz Ahmed","b Boult","1","3","0","0","33.33","Mumbai Indians","Dubai (DSC)","September 26","RCB won by 54 runs"),
    @("1st","Royal Challengers Bangalore","Shahbaz Ahmed","c KH Pandya b Jansen","1","2","0","0","50.00","Mumbai Indians","Chennai","April 09","RCB won by 2 wickets")
)

function Set-TextCell($cell, $text) {
    # Preserve numeric-looking strings ("13", "92.85", ...) as TEXT instead of
    # letting Excel auto-coerce them to numbers; non-numeric strings don't
    # need this so we leave their (default) style untouched.
    if ($text -match '^[0-9]+(\.[0-9]+)?$') {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $text
}

for ($c = 0; $c -lt $headers.Length; $c++) {
    Set-TextCell $ws.Cells.Item(1, $c + 1) $headers[$c]
}

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        Set-TextCell $ws.Cells.Item($r + 2, $c + 1) $row[$c]
    }
}
